# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) from SCD0271 to SCD0017
$ws.Name = "SCD0017"

# Update TC_ID values in column B (rows 2-4) from DGS-286 to SCD0017-001
$ws.Range("B2").Value = "SCD0017-001"
$ws.Range("B3").Value = "SCD0017-001"
$ws.Range("B4").Value = "SCD0017-001"

# Widen column B to fit the new, longer TC_ID text (closest width achievable
# to the target 13.140625 "best fit" character width)
$ws.Columns.Item(2).ColumnWidth = 12.25

# Update the active selection / view: the previous view had a frozen/scrolled
# left column at F1 with a selection at N2; the new view restores the left
# edge and selects B5 instead.
$ws.Range("B5").Select()
